# Generate Report for Handoff
# Replace references to the old localization-run UUID/hashes with the new
# ones, and bump the associated timestamps, mirroring a fresh handoff run.

$wb = $excel.ActiveWorkbook

$oldGuid = "6a362ab9-1a27-429b-a556-2d3c27da839f"
$newGuid = "b964ff0c-27b0-4326-8a2b-cb625594757d"

$oldHash = "95902da540d998d58804bd750a3cbcde72267fde"
$newHash = "2252a6a2323aa890a8ceef73037cfa85d9e21fa2"

$oldHoDate  = "2016-08-13 21:14:07"
$newHoDate  = "2016-08-13 21:14:37"

$oldZhDate  = "2016-08-13 21:13:56"
$newZhDate  = "2016-08-13 21:14:28"

# The external hyperlink targets (the .rels Target URLs) stay exactly as
# they were before the edit - only the visible cell text / hyperlink
# display text changes.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/164a01d4248262988f266de6c2aa1fd67b6d3119/e2e/$oldGuid.md"

function Update-HyperlinkDisplay {
    param($ws, $rangeAddr, $displayText)

    $r = $ws.Range($rangeAddr)
    # Drop the existing (external, read-only) hyperlink and recreate it
    # pointing at the same address but with the refreshed display text,
    # so the cell keeps exactly one hyperlink with the correct text.
    $r.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($r, $hyperlinkUrl, "", "", $displayText) | Out-Null
}

# NOTE: this runtime's PowerShell-style function calls only bind
# parameters positionally (`-name value` style named args do not bind),
# so every call below passes ws/rangeAddr/displayText in that order.

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# A2: File Name
$wsOverview.Range("A2").Value = "$newGuid.md"

# B2: Path And Name (also the hyperlink cell)
Update-HyperlinkDisplay $wsOverview "B2" "e2e\$newGuid.md"

# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = $newHoDate

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# A2: Source File Name (also the hyperlink cell)
Update-HyperlinkDisplay $wsZhCn "A2" "$newGuid.md"

# G2: Latest Handoff File
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"

# H2: Latest Handoff Datetime
$wsZhCn.Range("H2").Value = $newZhDate

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# A2: Source File Name (also the hyperlink cell)
Update-HyperlinkDisplay $wsDeDe "A2" "$newGuid.md"

# G2: Latest Handoff File
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"

# H2: Latest Handoff Datetime (shares the same text as Overview's G2)
$wsDeDe.Range("H2").Value = $newHoDate

Write-Host "Localization status report refreshed for handoff."
